$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Client -> Serveur, TL, Demande la liste des teams disponibles pour un match
$ws.Range("A13").Value = "Client"
$ws.Range("B13").Value = "Serveur"
$ws.Range("C13").Value = "TL"
$ws.Range("E13").Value = "Demande la liste des teams disponibles pour un match"

# Row 14: Serveur -> Client, TL, (voir code), Envoi la liste des équipes disponibles pour un match
$ws.Range("A14").Value = "Serveur"
$ws.Range("B14").Value = "Client"
$ws.Range("C14").Value = "TL"
$ws.Range("D14").Value = "(voir code)"
$ws.Range("E14").Value = "Envoi la liste des équipes disponibles pour un match"

# Row 15: Client -> Serveur, CM, Demande la création d'un match, matchName;idEquipe1;idEquipe2
$ws.Range("A15").Value = "Client"
$ws.Range("B15").Value = "Serveur"
$ws.Range("C15").Value = "CM"
$ws.Range("E15").Value = "Demande la création d'un match"
$ws.Range("D15").Value = "matchName;idEquipe1;idEquipe2"

# Row 16: Serveur -> Client, CF
$ws.Range("A16").Value = "Serveur"
$ws.Range("B16").Value = "Client"
$ws.Range("C16").Value = "CF"

# Row 17: Serveur -> Client, CO, Informe que le match a bien été créé.
$ws.Range("A17").Value = "Serveur"
$ws.Range("B17").Value = "Client"
$ws.Range("C17").Value = "CO"
$ws.Range("E17").Value = "Informe que le match a bien été créé."

# Back to row 16: description for CF
$ws.Range("E16").Value = "Informe que la création de match a échouée car les id d'équipe sont identiques."

# Row 18: Serveur -> Client, CN, Informe que le match ne peut pas etre cree car au moins une des deux équipes a déjà un match planifié.
$ws.Range("A18").Value = "Serveur"
$ws.Range("B18").Value = "Client"
$ws.Range("C18").Value = "CN"
$ws.Range("E18").Value = "Informe que le match ne peut pas etre cree car au moins une des deux équipes a déjà un match planifié."

# Row 19: Serveur -> Client, HW, Envoi le client vers l'écran d'attente d'un match.
$ws.Range("A19").Value = "Serveur"
$ws.Range("B19").Value = "Client"
$ws.Range("C19").Value = "HW"
$ws.Range("E19").Value = "Envoi le client vers l'écran d'attente d'un match."

# Row heights for the wrapped-text rows (matches taller rows in the saved sheet)
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 30
$ws.Rows.Item(18).RowHeight = 45

# Cursor/selection left where the author saved the file
$ws.Range("E19").Select()
